$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "57.132.92"
$ws.Range("E2").Value = "  -6.21%  "

$ws.Range("D3").Value = "2.891.88"
$ws.Range("E3").Value = "  -3.45%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "122.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.92%  "

$ws.Range("D8").Value = "2.890.40"
$ws.Range("E8").Value = "  -3.30%  "

$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.36%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.437"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.54%  "

$ws.Range("E13").Value = "  -5.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("E15").Value = "  +1.01%  "

$ws.Range("D16").Value = "3.367.59"
$ws.Range("E16").Value = "  -3.52%  "

$ws.Range("D17").Value = "2.888.78"
$ws.Range("E17").Value = "  -3.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.47%  "

$ws.Range("D19").Value = "57.182.39"
$ws.Range("E19").Value = "  -6.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "403.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "77.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("E28").Value = "  -1.78%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.61%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0990"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.911"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.04%  "

$ws.Range("E36").Value = "  -12.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.32%  "

$ws.Range("D39").Value = "0.0₃0621"
$ws.Range("E39").Value = "  -6.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0339"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.39%  "

$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("D42").Value = "2.627.47"
$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "358.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.52%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "119.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.229"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.79%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.63%  "
